$p = $ppt.ActivePresentation
$p.Slides.Item(7).MoveTo(6)
